$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix J2's cell border so it matches the rest of the J column (remove the
# unique thick-top border it had before, picking up the common format that
# J3:J23 already use). This mirrors the consolidation of cellXfs entries in
# the target file.
$ws.Range("J3").Copy() | Out-Null
$ws.Range("J2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Update the "Consultar Escenario" naming/document list in column J so that it
# matches the format described in Nomenclaturas.xlsx.
$ws.Range("J2").Value = "Plan de Iteración"
$ws.Range("J3").Value = "Informe de Revisión Técnica Formal"
$ws.Range("J4").Value = "Informe Final de SQA"
$ws.Range("J5").Value = "Consultar Escenario"
$ws.Range("J6").Value = "Asignar Escenario"
$ws.Range("J7").Value = "Adjuntar Documentos"
$ws.Range("J8").Value = "Comentar Escenario"

# Restore the selection to what the author left it at.
$ws.Range("J9").Select() | Out-Null

$excel.CutCopyMode = 0
